$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Draft 1" (list format) and "Draft 2" (paragraph format)
#    sections entirely -- everything from the start of the document up to
#    (but not including) the paragraph that begins the "Draft 3" section.
# ---------------------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Draft 3*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 1) {
    $startPara = $d.Paragraphs.Item(1)
    $endPara = $d.Paragraphs.Item($targetIndex - 1)
    $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $deleteRange.Delete()
}

# ---------------------------------------------------------------------------
# 2) In the "Legal Verifiers" bullet (Draft 3, Planned Employee Functions),
#    the description text was split across three runs ("... insurances and ",
#    "registrations", ", and run background checks"). Collapse it back into
#    a single run with the combined text, keeping the same formatting.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Legal Verifiers:*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $paraRange = $target.Range
    $marker = " Review and validate"
    $idx = $paraRange.Text.IndexOf($marker)
    if ($idx -ge 0) {
        $subStart = $paraRange.Start + $idx
        $subEnd = $paraRange.End - 1
        $finalText = " Review and validate customers' IDs or Passports, licenses, cars' insurances and registrations, and run background checks"

        # Assigning identical text to the existing Range.Text is a no-op (the
        # runs stay split), so first stamp a placeholder value to force a real
        # content change/run-merge, then set the final text.
        $subRange = $d.Range($subStart, $subEnd)
        $subRange.Text = $finalText + "#"
        $subRange2 = $d.Range($subStart, $subRange.End)
        $subRange2.Text = $finalText
    }
}
